# Update the "Statistical Description" sheet with refreshed statistics values.
# This mirrors a re-run of the underlying data analysis: the Mean/STD columns
# (C, D) were recomputed for every field, and several Quartile/Median values
# (F, G, H) shifted slightly as more data was incorporated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - co2
$ws.Range("C2").Value = 542.020306689403
$ws.Range("D2").Value = 130.1260856692631
$ws.Range("G2").Value = 495
$ws.Range("H2").Value = 603

# Row 3 - humidity
$ws.Range("C3").Value = 41.29301825412306
$ws.Range("D3").Value = 4.826245451811907
$ws.Range("F3").Value = 38.13
$ws.Range("G3").Value = 40.57
$ws.Range("H3").Value = 44.16

# Row 4 - pm25
$ws.Range("C4").Value = 1.417583508994922
$ws.Range("D4").Value = 1.991802455295512
$ws.Range("F4").Value = 0.52
$ws.Range("G4").Value = 1.01
$ws.Range("H4").Value = 1.82

# Row 5 - pressure
$ws.Range("C5").Value = 322.8550222308917
$ws.Range("D5").Value = 10.7791742118396
$ws.Range("F5").Value = 316.45
$ws.Range("G5").Value = 324.62
$ws.Range("H5").Value = 331.59

# Row 6 - temperature
$ws.Range("C6").Value = 20.77183497668683
$ws.Range("D6").Value = 2.582807601395347
$ws.Range("F6").Value = 19.4
$ws.Range("G6").Value = 20.85
$ws.Range("H6").Value = 22.26

# Row 7 - rssi
$ws.Range("C7").Value = -76.3014969513538
$ws.Range("D7").Value = 22.74087949442439

# Row 8 - snr
$ws.Range("C8").Value = 7.687664653310071
$ws.Range("D8").Value = 6.857906812619937

# Row 9 - SF
$ws.Range("C9").Value = 9.320027938058054
$ws.Range("D9").Value = 1.685034460598765

# Row 10 - frequency
$ws.Range("C10").Value = 867.8301732285452
$ws.Range("D10").Value = 0.4614636160985891

# Row 11 - toa
$ws.Range("C11").Value = 0.5549792477413591
$ws.Range("D11").Value = 0.5886060145151576

# Row 12 - distance
$ws.Range("C12").Value = 22.73730297062099
$ws.Range("D12").Value = 12.29191041355318

# Row 13 - c_walls
$ws.Range("C13").Value = 0.6738143629304758
$ws.Range("D13").Value = 0.7504754826094961

# Row 14 - w_walls
$ws.Range("C14").Value = 1.826669519201123
$ws.Range("D14").Value = 1.664155219542451

# Row 15 - exp_pl
$ws.Range("C15").Value = 93.70149695135368
$ws.Range("D15").Value = 22.74087949442439

# Row 16 - n_power
$ws.Range("C16").Value = -85.56402711020543
$ws.Range("D16").Value = 20.48085159011164
$ws.Range("H16").Value = -67.79706163635328

# Row 17 - esp
$ws.Range("C17").Value = -77.87636245689536
$ws.Range("D17").Value = 25.11899245527042
$ws.Range("F17").Value = -92.79009749652566
$ws.Range("G17").Value = -74.1773721860196
